$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# header1.xml == Headers.Item(2) (wdHeaderFooterFirstPage) -> BTec logo: image1.jpg -> image2.jpg
$hdrFirst = $sec.Headers.Item(2).Range.InlineShapes.Item(1)
$hdrFirst.Name = "image2.jpg"

# header2.xml == Headers.Item(1) (wdHeaderFooterPrimary) -> BTec logo: image1.jpg -> image2.jpg
$hdrPrimary = $sec.Headers.Item(1).Range.InlineShapes.Item(1)
$hdrPrimary.Name = "image2.jpg"

# footer1.xml == Footers.Item(2) (wdHeaderFooterFirstPage) -> Pearson logo: image2.png -> image1.png
$ftrFirst = $sec.Footers.Item(2).Range.InlineShapes.Item(1)
$ftrFirst.Name = "image1.png"

# footer2.xml == Footers.Item(1) (wdHeaderFooterPrimary) -> Pearson logo: image2.png -> image1.png
$ftrPrimary = $sec.Footers.Item(1).Range.InlineShapes.Item(1)
$ftrPrimary.Name = "image1.png"

Write-Host "Renamed header/footer inline shapes."
